$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Cronograma General")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy([System.Reflection.Missing]::Value, $last)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Quinta Semana"

# delete unwanted columns E:L (keep B:D); clear stray column A content (do not delete/shift)
$new.Columns("E:L").Delete()
$new.Columns("A:A").Clear()
# delete row1 (the spacer) entirely so header moves to row1
$new.Rows("1:1").Delete()
# now insert 3 blank rows above row1 so header moves to row4
$new.Rows("1:3").Insert()

Write-Host "dim:" $new.UsedRange.Address()
